$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2375.5
$ws.Range("J17").Value = 2375.5
$ws.Range("L17").Value = 7126.5
$ws.Range("N17").Value = -7462.5
$ws.Range("H19").Value = 1542.9445
$ws.Range("I19").Value = 666
$ws.Range("J19").Value = 1981.4166
$ws.Range("K19").Value = 666
$ws.Range("L19").Value = 1981.4166
$ws.Range("M19").Value = -491
$ws.Range("N19").Value = -2331.4166
$ws.Range("H38").Value = 6895.76
$ws.Range("I38").Value = 7024.9375
$ws.Range("K38").Value = 21074.8125
$ws.Range("M38").Value = -20702.8125
$ws.Range("H55").Value = 635.4545000000001
$ws.Range("I55").Value = 193.75
$ws.Range("J55").Value = 887.8570999999999
$ws.Range("K55").Value = 193.75
$ws.Range("L55").Value = 887.8570999999999
$ws.Range("M55").Value = 20.25
$ws.Range("N55").Value = -1315.8571
$ws.Range("H64").Value = 9099.700000000001
$ws.Range("J64").Value = 9529.058999999999
$ws.Range("L64").Value = 9529.058999999999
$ws.Range("N64").Value = -10025.059
$ws.Range("H67").Value = 9099.700000000001
$ws.Range("J67").Value = 9529.058999999999
$ws.Range("L67").Value = 9529.058999999999
$ws.Range("N67").Value = -11245.059
$ws.Range("H101").Value = 2783.8
$ws.Range("I101").Value = 979.75
$ws.Range("K101").Value = 2939.25
$ws.Range("M101").Value = -1317.25
$ws.Range("H106").Value = 4308.5454
$ws.Range("I106").Value = 3099.25
$ws.Range("K106").Value = 3099.25
$ws.Range("M106").Value = -2468.25
$ws.Range("H107").Value = 983.14703
$ws.Range("I107").Value = 936.0741
$ws.Range("K107").Value = 936.0741
$ws.Range("M107").Value = 983.9259
$ws.Range("H111").Value = 1853.0869
$ws.Range("I111").Value = 1485.4445
$ws.Range("J111").Value = 3176.6
$ws.Range("K111").Value = 4456.333500000001
$ws.Range("L111").Value = 9529.799999999999
$ws.Range("M111").Value = -1389.333500000001
$ws.Range("N111").Value = -15663.8
$ws.Range("H112").Value = 1586.8572
$ws.Range("I112").Value = 1052
$ws.Range("J112").Value = 2300
$ws.Range("K112").Value = 3156
$ws.Range("L112").Value = 6900
$ws.Range("M112").Value = -2048
$ws.Range("N112").Value = -9116
$ws.Range("H116").Value = 9356.6
$ws.Range("I116").Value = 8998.666999999999
$ws.Range("J116").Value = 9893.5
$ws.Range("K116").Value = 8998.666999999999
$ws.Range("L116").Value = 9893.5
$ws.Range("M116").Value = -5556.666999999999
$ws.Range("N116").Value = -16777.5
$ws.Range("H135").Value = 2495.3333
$ws.Range("I135").Value = 2327.111
$ws.Range("J135").Value = 3000
$ws.Range("K135").Value = 20943.999
$ws.Range("L135").Value = 27000
$ws.Range("M135").Value = -18408.999
$ws.Range("N135").Value = -32070
$ws.Range("H137").Value = 2298.1667
$ws.Range("I137").Value = 2573.75
$ws.Range("K137").Value = 7721.25
$ws.Range("M137").Value = -5171.25
$ws.Range("H138").Value = 8983.447
$ws.Range("J138").Value = 7944.643
$ws.Range("L138").Value = 23833.929
$ws.Range("N138").Value = -34113.929
$ws.Range("H141").Value = 4728.0625
$ws.Range("I141").Value = 4043.2666
$ws.Range("K141").Value = 12129.7998
$ws.Range("M141").Value = -6949.799800000001

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 11512.857
$ws.Range("I32").Value = 8692.352999999999
$ws.Range("J32").Value = 23500
$ws.Range("K32").Value = 8692.352999999999
$ws.Range("L32").Value = 23500
$ws.Range("M32").Value = -8405.352999999999
$ws.Range("N32").Value = -24074
$ws.Range("H45").Value = 9330.691999999999
$ws.Range("I45").Value = 10267.182
$ws.Range("K45").Value = 10267.182
$ws.Range("M45").Value = -9890.182000000001
$ws.Range("H61").Value = 6608.5
$ws.Range("I61").Value = 6320.1816
$ws.Range("J61").Value = 7665.6665
$ws.Range("K61").Value = 6320.1816
$ws.Range("L61").Value = 7665.6665
$ws.Range("M61").Value = -6108.1816
$ws.Range("N61").Value = -8089.6665
$ws.Range("H74").Value = 2674.1785
$ws.Range("I74").Value = 2590.68
$ws.Range("K74").Value = 2590.68
$ws.Range("M74").Value = -1716.68
$ws.Range("H76").Value = 23499
$ws.Range("J76").Value = 23499
$ws.Range("L76").Value = 23499
$ws.Range("N76").Value = -24175
$ws.Range("H77").Value = 2674.1785
$ws.Range("I77").Value = 2590.68
$ws.Range("K77").Value = 12953.4
$ws.Range("M77").Value = -8585.4
$ws.Range("H79").Value = 23499
$ws.Range("J79").Value = 23499
$ws.Range("L79").Value = 23499
$ws.Range("N79").Value = -25839
$ws.Range("H97").Value = 634
$ws.Range("I97").Value = 642.5
$ws.Range("K97").Value = 642.5
$ws.Range("M97").Value = -146.5
$ws.Range("H102").Value = 7494.0713
$ws.Range("I102").Value = 1731
$ws.Range("K102").Value = 1731
$ws.Range("M102").Value = -109
$ws.Range("H122").Value = 3500
$ws.Range("I122").Value = 4250
$ws.Range("K122").Value = 12750
$ws.Range("M122").Value = -10300
$ws.Range("H136").Value = 6608.5
$ws.Range("I136").Value = 6320.1816
$ws.Range("J136").Value = 7665.6665
$ws.Range("K136").Value = 18960.5448
$ws.Range("L136").Value = 22996.9995
$ws.Range("M136").Value = -16410.5448
$ws.Range("N136").Value = -28096.9995

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H51").Value = 39999.5
$ws.Range("J51").Value = 39999.5
$ws.Range("L51").Value = 39999.5
$ws.Range("N51").Value = -40981.5
$ws.Range("H54").Value = 23875
$ws.Range("I54").Value = 2750
$ws.Range("K54").Value = 2750
$ws.Range("M54").Value = -2266
$ws.Range("H94").Value = 1312.0769
$ws.Range("I94").Value = 1457.75
$ws.Range("J94").Value = 1079
$ws.Range("K94").Value = 1457.75
$ws.Range("L94").Value = 1079
$ws.Range("M94").Value = -1006.75
$ws.Range("N94").Value = -1981
$ws.Range("H99").Value = 4032.3076
$ws.Range("I99").Value = 2996.6667
$ws.Range("J99").Value = 4343
$ws.Range("K99").Value = 2996.6667
$ws.Range("L99").Value = 4343
$ws.Range("M99").Value = -1498.6667
$ws.Range("N99").Value = -7339
$ws.Range("H105").Value = 7324.2383
$ws.Range("I105").Value = 8153.933
$ws.Range("K105").Value = 8153.933
$ws.Range("M105").Value = -6406.933
$ws.Range("H107").Value = 4645.6875
$ws.Range("I107").Value = 836.625
$ws.Range("K107").Value = 836.625
$ws.Range("M107").Value = 1083.375
$ws.Range("H132").Value = 83150.5
$ws.Range("J132").Value = 83150.5
$ws.Range("L132").Value = 83150.5
$ws.Range("N132").Value = -93270.5
$ws.Range("H134").Value = 2591.5925
$ws.Range("I134").Value = 2284.476
$ws.Range("K134").Value = 6853.428
$ws.Range("M134").Value = -4318.428

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 8008.3076
$ws.Range("I7").Value = 263.33334
$ws.Range("J7").Value = 14646.857
$ws.Range("K7").Value = 263.33334
$ws.Range("L7").Value = 14646.857
$ws.Range("M7").Value = -150.33334
$ws.Range("N7").Value = -14872.857
$ws.Range("H22").Value = 2077.4443
$ws.Range("J22").Value = 2759.4
$ws.Range("L22").Value = 2759.4
$ws.Range("N22").Value = -3459.4
$ws.Range("H31").Value = 5497
$ws.Range("I31").Value = 2294
$ws.Range("K31").Value = 2294
$ws.Range("M31").Value = -1999
$ws.Range("H34").Value = 5497
$ws.Range("I34").Value = 2294
$ws.Range("K34").Value = 2294
$ws.Range("M34").Value = -2092
$ws.Range("H63").Value = 39999.89
$ws.Range("J63").Value = 39999.89
$ws.Range("L63").Value = 39999.89
$ws.Range("N63").Value = -41371.89
$ws.Range("H66").Value = 39999.89
$ws.Range("J66").Value = 39999.89
$ws.Range("L66").Value = 119999.67
$ws.Range("N66").Value = -126863.67
$ws.Range("H86").Value = 7629.727
$ws.Range("I86").Value = 7846.7144
$ws.Range("J86").Value = 7250
$ws.Range("K86").Value = 7846.7144
$ws.Range("L86").Value = 7250
$ws.Range("M86").Value = -6723.7144
$ws.Range("N86").Value = -9496
$ws.Range("H89").Value = 7629.727
$ws.Range("I89").Value = 7846.7144
$ws.Range("J89").Value = 7250
$ws.Range("K89").Value = 39233.572
$ws.Range("L89").Value = 36250
$ws.Range("M89").Value = -33617.572
$ws.Range("N89").Value = -47482
$ws.Range("H94").Value = 5082.4
$ws.Range("I94").Value = 4012
$ws.Range("J94").Value = 5350
$ws.Range("K94").Value = 4012
$ws.Range("L94").Value = 5350
$ws.Range("N94").Value = -6252
$ws.Range("H99").Value = 2165
$ws.Range("I99").Value = 1952.75
$ws.Range("K99").Value = 1952.75
$ws.Range("M99").Value = -454.75
$ws.Range("H107").Value = 1058.421
$ws.Range("I107").Value = 910.8889
$ws.Range("K107").Value = 910.8889
$ws.Range("M107").Value = 1009.1111
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("H126").Value = 2165
$ws.Range("I126").Value = 1952.75
$ws.Range("K126").Value = 5858.25
$ws.Range("M126").Value = -3388.25
$ws.Range("H132").Value = 1984.8889
$ws.Range("I132").Value = 1954
$ws.Range("J132").Value = 2000.3334
$ws.Range("K132").Value = 5862
$ws.Range("L132").Value = 6001.0002
$ws.Range("M132").Value = -3332
$ws.Range("N132").Value = -11061.0002
$ws.Range("H134").Value = 1715.2
$ws.Range("I134").Value = 1715.2
$ws.Range("K134").Value = 5145.6
$ws.Range("M134").Value = -2610.6
$ws.Range("M94").Value = -3561
$ws.Range("M122").ClearContents()
$ws.Range("N122").ClearContents()

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 202908.06
$ws.Range("I11").Value = 295737.28
$ws.Range("J11").Value = 371.54544
$ws.Range("K11").Value = 887211.8400000001
$ws.Range("L11").Value = 1114.63632
$ws.Range("M11").Value = -887071.8400000001
$ws.Range("N11").Value = -1394.63632
$ws.Range("H33").Value = 343.93332
$ws.Range("J33").Value = 387
$ws.Range("L33").Value = 2322
$ws.Range("N33").Value = -2888
$ws.Range("H57").Value = 5666.6665
$ws.Range("H64").Value = 1642.7778
$ws.Range("I64").Value = 714.1667
$ws.Range("J64").Value = 3500
$ws.Range("K64").Value = 2142.5001
$ws.Range("L64").Value = 10500
$ws.Range("M64").Value = -1872.5001
$ws.Range("N64").Value = -11040
$ws.Range("H67").Value = 1642.7778
$ws.Range("I67").Value = 714.1667
$ws.Range("J67").Value = 3500
$ws.Range("K67").Value = 2142.5001
$ws.Range("L67").Value = 10500
$ws.Range("M67").Value = -1206.5001
$ws.Range("N67").Value = -12372
$ws.Range("H121").Value = 4291.8184
$ws.Range("I121").Value = 852.5
$ws.Range("J121").Value = 6257.143
$ws.Range("K121").Value = 2557.5
$ws.Range("L121").Value = 18771.429
$ws.Range("M121").Value = -1247.5
$ws.Range("N121").Value = -21391.429
$ws.Range("H132").Value = 2999.1538
$ws.Range("I132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("H137").Value = 7723.16
$ws.Range("J137").Value = 7984.2
$ws.Range("L137").Value = 23952.6
$ws.Range("N137").Value = -34152.6
$ws.Range("M132").ClearContents()

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H17").Value = 65.125
$ws.Range("I17").Value = 51.8
$ws.Range("J17").Value = 87.333336
$ws.Range("K17").Value = 51.8
$ws.Range("L17").Value = 87.333336
$ws.Range("M17").Value = 116.2
$ws.Range("N17").Value = -423.333336
$ws.Range("H18").Value = 4004
$ws.Range("I18").Value = 5
$ws.Range("K18").Value = 5
$ws.Range("M18").Value = 288
$ws.Range("H40").Value = 12499.5
$ws.Range("I40").Value = 500
$ws.Range("K40").Value = 500
$ws.Range("I42").Value = 0
$ws.Range("K42").Value = 0
$ws.Range("H46").Value = 8749.75
$ws.Range("I46").Value = 8749.75
$ws.Range("K46").Value = 8749.75
$ws.Range("M46").Value = -8593.75
$ws.Range("H80").Value = 6485.1333
$ws.Range("I80").Value = 4826.6665
$ws.Range("J80").Value = 6899.75
$ws.Range("K80").Value = 4826.6665
$ws.Range("L80").Value = 6899.75
$ws.Range("M80").Value = -3828.6665
$ws.Range("N80").Value = -8895.75
$ws.Range("H82").Value = 75000
$ws.Range("J82").Value = 75000
$ws.Range("L82").Value = 75000
$ws.Range("N82").Value = -75766
$ws.Range("H83").Value = 6485.1333
$ws.Range("I83").Value = 4826.6665
$ws.Range("J83").Value = 6899.75
$ws.Range("K83").Value = 24133.3325
$ws.Range("L83").Value = 34498.75
$ws.Range("M83").Value = -19141.3325
$ws.Range("N83").Value = -44482.75
$ws.Range("H85").Value = 75000
$ws.Range("J85").Value = 75000
$ws.Range("L85").Value = 75000
$ws.Range("N85").Value = -77652
$ws.Range("H97").Value = 6147
$ws.Range("I97").Value = 7221
$ws.Range("J97").Value = 3999
$ws.Range("K97").Value = 7221
$ws.Range("L97").Value = 3999
$ws.Range("M97").Value = -6725
$ws.Range("N97").Value = -4991
$ws.Range("H113").Value = 15269.765
$ws.Range("I113").Value = 29897.25
$ws.Range("K113").Value = 29897.25
$ws.Range("M113").Value = -27727.25
$ws.Range("I115").Value = 0
$ws.Range("K115").Value = 0
$ws.Range("H132").Value = 1858.9117
$ws.Range("I132").Value = 1896.8
$ws.Range("J132").Value = 1574.75
$ws.Range("K132").Value = 5690.4
$ws.Range("L132").Value = 4724.25
$ws.Range("M132").Value = -3160.4
$ws.Range("N132").Value = -9784.25
$ws.Range("M40").Value = -349
$ws.Range("M42").ClearContents()
$ws.Range("M115").ClearContents()

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 3561.7856
$ws.Range("J46").Value = 3561.7856
$ws.Range("L46").Value = 3561.7856
$ws.Range("N46").Value = -3937.7856
$ws.Range("H62").Value = 400024200
$ws.Range("J62").Value = 400024200
$ws.Range("L62").Value = 400024200
$ws.Range("N62").Value = -400025448
$ws.Range("H65").Value = 400024200
$ws.Range("J65").Value = 400024200
$ws.Range("L65").Value = 1200072600
$ws.Range("N65").Value = -1200078840
$ws.Range("H82").Value = 4212
$ws.Range("I82").Value = 2246
$ws.Range("J82").Value = 4998.4
$ws.Range("K82").Value = 2246
$ws.Range("L82").Value = 4998.4
$ws.Range("M82").Value = -1885
$ws.Range("N82").Value = -5720.4
$ws.Range("H85").Value = 4212
$ws.Range("I85").Value = 2246
$ws.Range("J85").Value = 4998.4
$ws.Range("K85").Value = 2246
$ws.Range("L85").Value = 4998.4
$ws.Range("M85").Value = -998
$ws.Range("N85").Value = -7494.4
$ws.Range("H94").Value = 20000
$ws.Range("J94").Value = 20000
$ws.Range("L94").Value = 20000
$ws.Range("H100").Value = 8742
$ws.Range("I100").Value = 6880
$ws.Range("J100").Value = 10072
$ws.Range("K100").Value = 6880
$ws.Range("L100").Value = 10072
$ws.Range("M100").Value = -6339
$ws.Range("N100").Value = -11154
$ws.Range("H120").Value = 127566
$ws.Range("J120").Value = 127566
$ws.Range("L120").Value = 127566
$ws.Range("N120").Value = -137242
$ws.Range("H132").Value = 4808.25
$ws.Range("I132").Value = 2633.3333
$ws.Range("J132").Value = 5533.222
$ws.Range("K132").Value = 7899.999899999999
$ws.Range("L132").Value = 16599.666
$ws.Range("M132").Value = -5369.999899999999
$ws.Range("N132").Value = -21659.666
$ws.Range("H136").Value = 17096.584
$ws.Range("I136").Value = 2986.625
$ws.Range("J136").Value = 129976.25
$ws.Range("K136").Value = 8959.875
$ws.Range("L136").Value = 389928.75
$ws.Range("M136").Value = -6409.875
$ws.Range("N136").Value = -395028.75
$ws.Range("N94").Value = -21352

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H58").Value = 0
$ws.Range("J58").Value = 0
$ws.Range("L58").Value = 0
$ws.Range("H96").Value = 3121.2942
$ws.Range("I96").Value = 1206.2
$ws.Range("K96").Value = 1206.2
$ws.Range("M96").Value = 166.8
$ws.Range("H100").Value = 1313.8572
$ws.Range("I100").Value = 539.4
$ws.Range("K100").Value = 1078.8
$ws.Range("M100").Value = -537.8
$ws.Range("H101").Value = 14372.875
$ws.Range("J101").Value = 14372.875
$ws.Range("L101").Value = 14372.875
$ws.Range("N101").Value = -20862.875
$ws.Range("H107").Value = 482.93332
$ws.Range("I107").Value = 144.4
$ws.Range("J107").Value = 652.2
$ws.Range("K107").Value = 433.2
$ws.Range("L107").Value = 1956.6
$ws.Range("M107").Value = 1486.8
$ws.Range("N107").Value = -5796.6
$ws.Range("H113").Value = 592.4643
$ws.Range("I113").Value = 549.4706
$ws.Range("J113").Value = 658.9091
$ws.Range("K113").Value = 1648.4118
$ws.Range("L113").Value = 1976.7273
$ws.Range("M113").Value = 521.5882000000001
$ws.Range("N113").Value = -6316.7273
$ws.Range("N58").ClearContents()
